$d = $word.ActiveDocument

# Update the date heading (unique text in the document, safe to use Find/Replace).
$d.Content.Find.Execute("2023-10-22 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-10-23 Monday", 2)

# Update the practice-problem table. Addressing cells positionally (row, column)
# and assigning directly to Cell.Range.Text avoids ambiguity since some old
# values (e.g. "67÷2=") repeat in the grid but map to different new values in
# different cells; Find/Execute on a sub-range is not reliably scoped to that
# range in this environment, so direct Range.Text assignment (which *is*
# properly scoped) is used instead.
$t = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)
$newVals = @(
    @("86÷4=", "55÷8=", "87÷9=", "55÷6=", "53÷7="),
    @("47÷9=", "95÷7=", "22÷4=", "73÷2=", "65÷5="),
    @("52÷7=", "17÷7=", "85÷8=", "80÷2=", "40÷4="),
    @("29÷2=", "15÷4=", "98÷3=", "24÷7=", "12÷8="),
    @("11÷3=", "77÷9=", "31÷2=", "36÷7=", "83÷7=")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowIdx = $rows[$i]
    for ($colIdx = 1; $colIdx -le 5; $colIdx++) {
        $cell = $t.Cell($rowIdx, $colIdx)
        $cell.Range.Text = $newVals[$i][$colIdx - 1]
    }
}
